# Applies the "APELLIDO ESPOSO" column insert edit described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (pushes D..O -> E..P).
$ws.Range("D1").EntireColumn.Insert()

# New column D inherits column C's header text context: set the new header values.
$ws.Range("C5").Value = "APELLIDO ESPOSO"
$ws.Range("D5").Value = "CI"
$ws.Range("D5").Style = "Normal"

$ws.Range("H5").Value = "USUARIO"

# Restore view: active cell C7, top-left A1
$ws.Range("C7").Select()
